# The deck originally ships with the "Integral" (Red Violet) design applied
# to the slide master / notes master, and keeps an unused "Office Theme"
# theme part around (only reachable through the notes master's own theme
# relationship). The edit swaps the two: the presentation's applied design
# goes back to the default "Office Theme" colour palette, while the
# secondary theme part is reset to hold the "Integral" / "Red Violet"
# palette that used to be active.
#
# PowerPoint's theme colours are exposed as OLE (BGR-packed) RGB integers
# through ThemeColorScheme.Colors(index).RGB, so convert each target hex
# triplet (as it appears in the OOXML <a:srgbClr val="RRGGBB"/>) before
# assigning it.

function ConvertTo-OleColor([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return ($b * 65536) + ($g * 256) + $r
}

# ThemeColorScheme.Colors() index order: dk1, lt1, dk2, lt2, accent1-6,
# hlink, folHlink (MsoThemeColorSchemeIndex 1-12).
$officeThemeColors = @(
    "000000", # dk1
    "FFFFFF", # lt1
    "44546A", # dk2
    "E7E6E6", # lt2
    "5B9BD5", # accent1
    "ED7D31", # accent2
    "A5A5A5", # accent3
    "FFC000", # accent4
    "4472C4", # accent5
    "70AD47", # accent6
    "0563C1", # hlink
    "954F72"  # folHlink
)

$integralThemeColors = @(
    "000000", # dk1
    "FFFFFF", # lt1
    "454551", # dk2
    "D8D9DC", # lt2
    "E32D91", # accent1
    "C830CC", # accent2
    "4EA6DC", # accent3
    "4775E7", # accent4
    "8971E1", # accent5
    "D54773", # accent6
    "6B9F25", # hlink
    "8C8C8C"  # folHlink
)

$p = $ppt.ActivePresentation

# Reset the notes master's theme colours to the "Integral" palette first...
$notesMaster = $p.NotesMaster
$notesScheme = $notesMaster.Theme.ThemeColorScheme
for ($i = 1; $i -le 12; $i++) {
    $notesScheme.Colors($i).RGB = ConvertTo-OleColor $integralThemeColors[$i - 1]
}

# ...then apply the "Office Theme" palette to the slide master (the design
# actually shown on the slides), so the presentation's visible theme ends up
# back on the default Office colours.
$slideMaster = $p.SlideMaster
$masterScheme = $slideMaster.Theme.ThemeColorScheme
for ($i = 1; $i -le 12; $i++) {
    $masterScheme.Colors($i).RGB = ConvertTo-OleColor $officeThemeColors[$i - 1]
}
